$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.759.26"
$ws.Range("E2").Value = "  +4.18%  "

$ws.Range("D3").Value = "2.640.19"
$ws.Range("E3").Value = "  +2.94%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'606.71"
$ws.Range("E5").Value = "  +0.90%  "

$ws.Range("D6").Value = "'180.03"
$ws.Range("E6").Value = "  +0.48%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("E8").Value = "  +1.79%  "

$ws.Range("E9").Value = "  +11.56%  "

$ws.Range("D10").Value = "2.638.83"
$ws.Range("E10").Value = "  +2.95%  "

$ws.Range("E12").Value = "  +2.82%  "

$ws.Range("E13").Value = "  +0.71%  "

$ws.Range("E14").Value = "  +5.22%  "

$ws.Range("D15").Value = "3.106.31"

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "72.541.06"
$ws.Range("E16").Value = "  +3.99%  "

$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").Value = "'26.87"
$ws.Range("E17").Value = "  +2.12%  "

$ws.Range("D18").Value = "2.635.04"
$ws.Range("E18").Value = "  +2.73%  "

$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'384.28"
$ws.Range("E19").Value = "  +5.05%  "

$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'11.72"
$ws.Range("E20").Value = "  +4.96%  "

$ws.Range("D21").Value = "'7.98"
$ws.Range("E21").Value = "  +3.98%  "

$ws.Range("E22").Value = "  +1.51%  "

$ws.Range("D23").Value = "'2.05"
$ws.Range("E23").Value = "  +16.66%  "

$ws.Range("D24").Value = "'74.13"
$ws.Range("E24").Value = "  +4.83%  "

$ws.Range("E25").Value = "  +2.94%  "

$ws.Range("E26").Value = "  +0.19%  "

$ws.Range("D27").Value = "'10.03"
$ws.Range("E27").Value = "  +8.36%  "

$ws.Range("D28").Value = "2.774.38"
$ws.Range("E28").Value = "  +3.11%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("D30").Value = "0.0₃0964"
$ws.Range("E30").Value = "  +3.78%  "

$ws.Range("D31").Value = "'8.13"
$ws.Range("E31").Value = "  +4.17%  "

$ws.Range("D32").Value = "'520.32"
$ws.Range("E32").Value = "  -0.08%  "

$ws.Range("E33").Value = "  +4.06%  "

$ws.Range("E34").Value = "  +2.02%  "

$ws.Range("E35").Value = "  -0.17%  "

$ws.Range("D36").Value = "'165.22"
$ws.Range("E36").Value = "  +1.29%  "

$ws.Range("E37").Value = "  +2.54%  "

$ws.Range("E38").Value = "  +4.55%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.112"
$ws.Range("E39").Value = "  -6.62%  "

$ws.Range("B40").Value = "WhiteBITCoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D40").Value = "'19.09"
$ws.Range("E40").Value = "  +0.90%  "

$ws.Range("E41").Value = "  +5.50%  "

$ws.Range("E42").Value = "  +4.35%  "

$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("E44").Value = "  +4.65%  "

$ws.Range("D45").Value = "'0.336"
$ws.Range("E45").Value = "  +2.83%  "

$ws.Range("D46").Value = "'39.47"
$ws.Range("E46").Value = "  +0.96%  "

$ws.Range("D47").Value = "'150.92"
$ws.Range("E47").Value = "  -1.53%  "

$ws.Range("E48").Value = "  +2.42%  "

$ws.Range("E49").Value = "  +4.38%  "

$ws.Range("E50").Value = "  +4.32%  "

$ws.Range("E51").Value = "  +2.95%  "
